$wb = $excel.ActiveWorkbook

# The new "UI Test" sheet mirrors the layout/styles of "Session Test" (same
# headers, column widths, row heights and cell styles), so create it by
# copying that sheet to the end of the workbook, then rename it and replace
# row 4 with the new UI test case content.
$sessionSheet = $wb.Worksheets.Item("Session Test")
$sessionSheet.Copy($null, $sessionSheet)

$newSheet = $wb.ActiveSheet
$newSheet.Name = "UI Test"

# Fill row 4 in the same order a person tabbing across the form would use
# (ID columns first, then the descriptive columns) so newly interned shared
# strings land in the same order as the authored workbook.
$newSheet.Range("C4").Value = "UI-01"
$newSheet.Range("E4").Value = "UI-TC-01"
$newSheet.Range("G4").Value = "UI-TS-01"
$newSheet.Range("F4").Value = "UI Module"
$newSheet.Range("D4").Value = "Verify error message is clearly visible to the user"
$newSheet.Range("H4").Value = "Verify error message visibility"
$newSheet.Range("I4").Value = "User is on login or registration page"
$newSheet.Range("J4").Value = "1. Submit form with invalid data `n2. Observe displayed error message"
$newSheet.Range("K4").Value = "Invalid email/password"
$newSheet.Range("L4").Value = "Error message is clearly visible, readable, and positioned near related field"
$newSheet.Range("M4").Value = "Medium"

# Match the saved selection on the new sheet.
$newSheet.Range("M4").Select()
